$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cavity detection algo")

# Add new row 7 data
$ws.Range("B7").Value = "center diff is 0"
$ws.Range("F7").Value = 65
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 2

# Add a comment to F7
$c = $ws.Range("F7").AddComment("Author:" + [char]10 + [char]10 + "False negatives:" + [char]10 + "-ADB1_51")

# Update selection
$ws.Range("K14").Select()
